$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the data range so numeric-looking strings
# (e.g. "1.00", "592.78") are not silently coerced to numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '60.999.54'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '2.920.13'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '592.78'
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '146.11'
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '0.507'
$ws.Range("E8").Value = '  +0.39%  '
$ws.Range("B9").Value = 'Toncoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D9").Value = '6.92'
$ws.Range("E9").Value = '  +1.80%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.144'
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").Value = '0.439'
$ws.Range("E11").Value = '  -2.25%  '
$ws.Range("B12").Value = 'ShibaInu'
$ws.Range("C12").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D12").Value = '0.0000226'
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = '33.63'
$ws.Range("E13").Value = '  -0.99%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '0.127'
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.406.03'
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '60.981.70'
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = '6.71'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.921.35'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '432.86'
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '13.39'
$ws.Range("E20").Value = '  -2.17%  '
$ws.Range("B21").Value = 'Polygon'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D21").Value = '0.679'
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '7.03'
$ws.Range("E22").Value = '  -1.79%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '81.81'
$ws.Range("E23").Value = '  +1.44%  '
$ws.Range("B24").Value = 'RenderToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D24").Value = '11.04'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("B25").Value = 'Fetch.AI'
$ws.Range("C25").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D25").Value = '2.20'
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '11.85'
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = '2.30'
$ws.Range("E28").Value = '  +5.75%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '2.61'
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = '7.03'
$ws.Range("E30").Value = '  -3.56%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '26.58'
$ws.Range("E31").Value = '  -0.35%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.107'
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("B34").Value = 'PEPE'
$ws.Range("C34").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D34").Value = '0.0₃0863'
$ws.Range("E34").Value = '  -1.74%  '
$ws.Range("B35").Value = 'Mantle'
$ws.Range("C35").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D35").Value = '1.01'
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").Value = '5.63'
$ws.Range("E36").Value = '  -0.62%  '
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").Value = '3.02'
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '0.123'
$ws.Range("E38").Value = '  -2.84%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '1.99'
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("B40").Value = 'Cosmos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D40").Value = '8.56'
$ws.Range("E40").Value = '  -1.60%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").Value = '0.288'
$ws.Range("E41").Value = '  -2.08%  '
$ws.Range("B42").Value = 'Arweave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D42").Value = '41.04'
$ws.Range("E42").Value = '  -2.08%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '375.89'
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0344'
$ws.Range("E44").Value = '  -1.35%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.701.56'
$ws.Range("E45").Value = '  +1.01%  '
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").Value = '132.92'
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = '23.84'
$ws.Range("E48").Value = '  -3.20%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '0.106'
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '2.02'
$ws.Range("E50").Value = '  -2.53%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.125'
$ws.Range("E51").Value = '  +0.33%  '

# Restore default (unstyled) formatting now that the text values are locked in,
# matching the original workbook where these cells carry no explicit style.
$dataRange.ClearFormats()
